# Daily attendance processing - 2026-01-23 11:10:06
# Normalizes the "Recorded By" column (G) so that the user's email address
# is listed before the literal "System" token, e.g.
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

$target = "System, dnasr281@gmail.com"
$replacement = "dnasr281@gmail.com, System"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $target) {
        $cell.Value = $replacement
    }
}
